$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$regcntr = 10002
$machine = 10032
$startDevice = 3000176

for ($i = 0; $i -lt 5; $i++) {
    $row = 157 + $i
    $ws.Cells.Item($row, 1).Value = $regcntr
    $ws.Cells.Item($row, 2).Value = $machine
    $ws.Cells.Item($row, 3).Value = $startDevice + $i
    $ws.Cells.Item($row, 4).Value = "eng"
    $ws.Cells.Item($row, 5).Value = $true
    $ws.Cells.Item($row, 6).Value = "superadmin"
    $ws.Cells.Item($row, 7).Value = "now()"
}

$ws.Range("E157").Select()
$excel.ActiveWindow.ScrollRow = 150
$excel.ActiveWindow.ScrollColumn = 1
